# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# header style already used on the existing header row (e.g. H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from an existing header cell
# so the new header cells reuse the same style record as B1:H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-64
$data = @{
    2  = @(13, 13)
    3  = @(9, 9)
    4  = @(7, 8)
    5  = @(10, 10)
    6  = @(7, 7)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(8, 8)
    11 = @(6, 6)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(7, 8)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(9, 9)
    22 = @(7, 7)
    23 = @(9, 9)
    24 = @(7, 7)
    25 = @(9, 9)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(9, 9)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(7, 7)
    35 = @(7, 7)
    36 = @(1, 2)
    37 = @(8, 8)
    38 = @(7, 7)
    39 = @(8, 8)
    40 = @(7, 7)
    41 = @(9, 9)
    42 = @(7, 7)
    43 = @(8, 8)
    44 = @(8, 8)
    45 = @(7, 7)
    46 = @(6, 6)
    47 = @(7, 7)
    48 = @(8, 8)
    49 = @(6, 6)
    50 = @(7, 8)
    51 = @(6, 6)
    52 = @(6, 6)
    53 = @(8, 8)
    54 = @(8, 8)
    55 = @(6, 6)
    56 = @(5, 6)
    57 = @(8, 9)
    58 = @(7, 7)
    59 = @(5, 6)
    60 = @(8, 8)
    61 = @(8, 8)
    62 = @(9, 9)
    63 = @(6, 6)
    64 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
